# Update benchmark: 2025-11-28 06:41:25 UTC
# Applies the per-cell text changes described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SANS OYUNLARI
$ws.Range("D2").Value = "23,81 TL - 23,81 TL"

# Row 3 - HESAPTAN EFT - Sube
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DUZENLI EFT
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8 - HESAPTAN HAVALE - Sube
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11 - DUZENLI HAVALE
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12 - GIDEN SWIFT
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"

# Row 14 - GIDEN SWIFT - Mobil
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"

# Row 15 - CEK TAHSILI BASKA BANKA
$ws.Range("D15").Value = "%0,8 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 3.047,62 TL"

# Row 17 - AYNI SUBE CEK TAHSILATI
$ws.Range("D17").Value = "%0,8 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 3.047,62 TL"

# Row 20 - CEK IADE
$ws.Range("D20").Value = "285,72 TL"

# Row 21 - BLOKE CEK DUZENLEME
$ws.Range("D21").Value = "%0,5 Asgari Tutar: 428,58 TL Azami Tutar: 428,58 TL / 5.523,81 TL"

# Row 22 - YP CEK TAKASA GONDERME
$ws.Range("D22").Value = "%1 Asgari Tutar: 285,72 TL Azami Tutar: 285,72 TL / 6.857,15 TL"

# Row 23 - CEK KARNESI SAYFA UCRETI
$ws.Range("D23").Value = "64,77 TL"

# Row 24 - SENET TAHSILE ALMA
$ws.Range("D24").Value = "476,2 TL"

# Row 25 - MUAMELESIZ SENET IADESI
$ws.Range("D25").Value = "428,58 TL"
